$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8192.875
$ws.Range("I40").Value = 8866
$ws.Range("J40").Value = 7789
$ws.Range("K40").Value = 8866
$ws.Range("L40").Value = 7789
$ws.Range("M40").Value = -8691
$ws.Range("N40").Value = -8139
$ws.Range("H106").Value = 1252.1818
$ws.Range("I106").Value = 1252.1818
$ws.Range("K106").Value = 1252.1818
$ws.Range("M106").Value = -621.1818000000001
$ws.Range("H113").Value = 2308.55
$ws.Range("I113").Value = 2215.2354
$ws.Range("K113").Value = 2215.2354
$ws.Range("M113").Value = 1038.7646
$ws.Range("H127").Value = 1497.7333
$ws.Range("I127").Value = 1497.7333
$ws.Range("K127").Value = 4493.199900000001
$ws.Range("M127").Value = 466.8000999999995
$ws.Range("H129").Value = 1771.0392
$ws.Range("J129").Value = 1837.9773
$ws.Range("L129").Value = 5513.9319
$ws.Range("N129").Value = -15513.9319
$ws.Range("H131").Value = 3182.3333
$ws.Range("I131").Value = 2773.75
$ws.Range("J131").Value = 3999.5
$ws.Range("K131").Value = 8321.25
$ws.Range("L131").Value = 11998.5
$ws.Range("M131").Value = -3281.25
$ws.Range("N131").Value = -22078.5
$ws.Range("H133").Value = 99993.664
$ws.Range("J133").Value = 99993.664
$ws.Range("L133").Value = 99993.664
$ws.Range("N133").Value = -110113.664
$ws.Range("H135").Value = 687.5897
$ws.Range("I135").Value = 478.75
$ws.Range("K135").Value = 4308.75
$ws.Range("M135").Value = -1773.75
$ws.Range("H137").Value = 6260.625
$ws.Range("I137").Value = 3347.5
$ws.Range("K137").Value = 10042.5
$ws.Range("M137").Value = -7492.5
$ws.Range("H138").Value = 2311.8333
$ws.Range("I138").Value = 2107.6
$ws.Range("K138").Value = 6322.799999999999
$ws.Range("M138").Value = -1182.799999999999
$ws.Range("H141").Value = 2252.2144
$ws.Range("I141").Value = 2117.111
$ws.Range("K141").Value = 6351.333
$ws.Range("M141").Value = -1171.333

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3388.7144
$ws.Range("I2").Value = 1898.6666
$ws.Range("J2").Value = 4506.25
$ws.Range("K2").Value = 1898.6666
$ws.Range("L2").Value = 4506.25
$ws.Range("M2").Value = -1785.6666
$ws.Range("N2").Value = -4732.25
$ws.Range("H32").Value = 1603.6522
$ws.Range("I32").Value = 1603.6522
$ws.Range("K32").Value = 1603.6522
$ws.Range("M32").Value = -1316.6522
$ws.Range("H94").Value = 94888.25
$ws.Range("J94").Value = 94888.25
$ws.Range("L94").Value = 94888.25
$ws.Range("N94").Value = -96690.25
$ws.Range("H116").Value = 3388.7144
$ws.Range("I116").Value = 1898.6666
$ws.Range("J116").Value = 4506.25
$ws.Range("K116").Value = 1898.6666
$ws.Range("L116").Value = 4506.25
$ws.Range("M116").Value = 395.3334
$ws.Range("N116").Value = -9094.25
$ws.Range("H132").Value = 83335830
$ws.Range("I132").Value = 2994.5
$ws.Range("J132").Value = 250001500
$ws.Range("K132").Value = 8983.5
$ws.Range("L132").Value = 750004500
$ws.Range("M132").Value = -6453.5
$ws.Range("N132").Value = -750009560

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3388.7144
$ws.Range("I3").Value = 1898.6666
$ws.Range("J3").Value = 4506.25
$ws.Range("K3").Value = 1898.6666
$ws.Range("L3").Value = 4506.25
$ws.Range("M3").Value = -1784.6666
$ws.Range("N3").Value = -4734.25
$ws.Range("H94").Value = 2460.3333
$ws.Range("I94").Value = 2209.25
$ws.Range("J94").Value = 2962.5
$ws.Range("K94").Value = 2209.25
$ws.Range("L94").Value = 2962.5
$ws.Range("M94").Value = -1758.25
$ws.Range("N94").Value = -3864.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1666.3
$ws.Range("I2").Value = 1851.625
$ws.Range("K2").Value = 1851.625
$ws.Range("M2").Value = -1738.625
$ws.Range("H18").Value = 34876.5
$ws.Range("J18").Value = 34876.5
$ws.Range("L18").Value = 34876.5
$ws.Range("N18").Value = -35336.5
$ws.Range("H22").Value = 693.4
$ws.Range("I22").Value = 915.6667
$ws.Range("J22").Value = 360
$ws.Range("K22").Value = 915.6667
$ws.Range("L22").Value = 360
$ws.Range("M22").Value = -565.6667
$ws.Range("N22").Value = -1060
$ws.Range("H31").Value = 2277.5833
$ws.Range("I31").Value = 2277.5833
$ws.Range("K31").Value = 2277.5833
$ws.Range("M31").Value = -1982.5833
$ws.Range("H34").Value = 2277.5833
$ws.Range("I34").Value = 2277.5833
$ws.Range("K34").Value = 2277.5833
$ws.Range("M34").Value = -2075.5833
$ws.Range("H62").Value = 10288.6
$ws.Range("I62").Value = 2860.75
$ws.Range("J62").Value = 40000
$ws.Range("K62").Value = 2860.75
$ws.Range("L62").Value = 40000
$ws.Range("M62").Value = -2236.75
$ws.Range("N62").Value = -41248
$ws.Range("H65").Value = 10288.6
$ws.Range("I65").Value = 2860.75
$ws.Range("J65").Value = 40000
$ws.Range("K65").Value = 14303.75
$ws.Range("L65").Value = 200000
$ws.Range("M65").Value = -11183.75
$ws.Range("N65").Value = -206240
$ws.Range("H132").Value = 8476.8125
$ws.Range("I132").Value = 9384.909
$ws.Range("J132").Value = 6479
$ws.Range("K132").Value = 28154.727
$ws.Range("L132").Value = 19437
$ws.Range("M132").Value = -25624.727
$ws.Range("N132").Value = -24497
$ws.Range("H134").Value = 7696986.5
$ws.Range("I134").Value = 3966.5
$ws.Range("K134").Value = 11899.5
$ws.Range("M134").Value = -9364.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 208.08
$ws.Range("I12").Value = 23.6
$ws.Range("J12").Value = 331.06668
$ws.Range("K12").Value = 70.80000000000001
$ws.Range("L12").Value = 993.2000400000001
$ws.Range("M12").Value = 102.2
$ws.Range("N12").Value = -1339.20004
$ws.Range("H98").Value = 403
$ws.Range("J98").Value = 545.6667
$ws.Range("L98").Value = 1637.0001
$ws.Range("N98").Value = -4633.0001
$ws.Range("H107").Value = 749.7059
$ws.Range("J107").Value = 904.9167
$ws.Range("L107").Value = 2714.7501
$ws.Range("N107").Value = -6554.7501
$ws.Range("H123").Value = 5749
$ws.Range("H132").Value = 1487.7273
$ws.Range("I132").Value = 627.25
$ws.Range("J132").Value = 1979.4286
$ws.Range("K132").Value = 5645.25
$ws.Range("L132").Value = 17814.8574
$ws.Range("M132").Value = -3115.25
$ws.Range("N132").Value = -22874.8574
$ws.Range("H139").Value = 3247.0605
$ws.Range("I139").Value = 2693.1333
$ws.Range("K139").Value = 8079.3999
$ws.Range("M139").Value = -2939.3999
$ws.Range("H140").Value = 2693.5715
$ws.Range("I140").Value = 2670
$ws.Range("K140").Value = 8010
$ws.Range("M140").Value = -2830

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5463869
$ws.Range("I11").Value = 6548642
$ws.Range("J11").Value = 40004
$ws.Range("K11").Value = 6548642
$ws.Range("L11").Value = 40004
$ws.Range("M11").Value = -6548503
$ws.Range("N11").Value = -40282
$ws.Range("H12").Value = 1200
$ws.Range("J12").Value = 1200
$ws.Range("L12").Value = 1200
$ws.Range("N12").Value = -1480
$ws.Range("H14").Value = 8333983.5
$ws.Range("I14").Value = 8333983.5
$ws.Range("K14").Value = 8333983.5
$ws.Range("M14").Value = -8333815.5
$ws.Range("H107").Value = 730.2
$ws.Range("I107").Value = 578.4375
$ws.Range("K107").Value = 578.4375
$ws.Range("M107").Value = 1341.5625
$ws.Range("H126").Value = 4270.524
$ws.Range("I126").Value = 6024
$ws.Range("K126").Value = 18072
$ws.Range("M126").Value = -15602
$ws.Range("H132").Value = 1933
$ws.Range("I132").Value = 1933
$ws.Range("K132").Value = 5799
$ws.Range("M132").Value = -3269

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 5857.8887
$ws.Range("J9").Value = 12490.75
$ws.Range("L9").Value = 12490.75
$ws.Range("N9").Value = -12938.75
$ws.Range("H22").Value = 1312.2142
$ws.Range("J22").Value = 1545.9231
$ws.Range("L22").Value = 1545.9231
$ws.Range("N22").Value = -2135.9231
$ws.Range("H27").Value = 1312.2142
$ws.Range("J27").Value = 1545.9231
$ws.Range("L27").Value = 1545.9231
$ws.Range("N27").Value = -1759.9231
$ws.Range("H61").Value = 1080.1428
$ws.Range("I61").Value = 977.1111
$ws.Range("K61").Value = 977.1111
$ws.Range("M61").Value = -775.1111
$ws.Range("H113").Value = 1080.1428
$ws.Range("I113").Value = 977.1111
$ws.Range("K113").Value = 977.1111
$ws.Range("M113").Value = 1192.8889
$ws.Range("H122").Value = 3513.1333
$ws.Range("J122").Value = 3683.889
$ws.Range("L122").Value = 11051.667
$ws.Range("N122").Value = -15951.667
$ws.Range("H132").Value = 2497.6365
$ws.Range("J132").Value = 3494.5
$ws.Range("L132").Value = 10483.5
$ws.Range("N132").Value = -15543.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 7292.6
$ws.Range("J8").Value = 9750
$ws.Range("L8").Value = 9750
$ws.Range("N8").Value = -10030
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 5000
$ws.Range("K17").Value = 5000
$ws.Range("M17").Value = -4828
$ws.Range("H132").Value = 1399.5454
$ws.Range("I132").Value = 1399.5454
$ws.Range("K132").Value = 4198.6362
$ws.Range("M132").Value = -1668.6362
